# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled the same as the other header cells (copy format from G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save values for rows 2-54 (1 = saved, 0 = not saved)
$saveValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 1
    53 = 0
    54 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
